$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumen")

# Update row 13 ("tierras_aridos") values per corrected figures sent by email
$ws.Range("H13").Value = 38
$ws.Range("J13").Value = 286
$ws.Range("K13").Value = 17.2
$ws.Range("L13").Value = 2
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 2
$ws.Range("Q13").Value = 8
